$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @{
  2  = @(44497, 400, 11500, 12000, 11750, 5875)
  3  = @(44462, 140, 13000, 14000, 13500, 6750)
  4  = @(44495, 300, 11000, 12000, 11500, 5750)
  5  = @(44463, 100, 13000, 14000, 13500, 6750)
  6  = @(44452, 200, 13000, 14000, 13500, 6750)
  7  = @(44498, 240, 11000, 11500, 11250, 5625)
  8  = @(44448, 100, 14000, 15000, 14500, 7250)
  9  = @(44489, 400, 11500, 12000, 11750, 5875)
  10 = @(44455, 160, 13000, 14000, 13500, 6750)
  11 = @(44446, 300, 14000, 15000, 14500, 7250)
  12 = @(44468, 300, 13000, 14000, 13500, 6750)
  13 = @(44459, 160, 13000, 14000, 13500, 6750)
  14 = @(44454, 300, 13000, 14000, 13500, 6750)
  15 = @(44490, 160, 11500, 12000, 11750, 5875)
  16 = @(44466, 160, 13500, 14000, 13750, 6875)
  17 = @(44494, 200, 11500, 12000, 11750, 5875)
  18 = @(44445, 160, 14000, 15000, 14500, 7250)
  19 = @(44491, 200, 11500, 12000, 11750, 5875)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals[0]
  $ws.Range("M$r").Value = $vals[1]
  $ws.Range("N$r").Value = $vals[2]
  $ws.Range("O$r").Value = $vals[3]
  $ws.Range("P$r").Value = $vals[4]
  $ws.Range("S$r").Value = $vals[5]
}
